# "adding the sudan data plus cleaning code and dictionnary"
#
# - add a new "sudan" worksheet (after "nigeria") with the dictionary rows
#   for the sudan indicators (new shared-string entries get created
#   automatically as the cell values are written)
# - restore/clean up the selections on the existing sheets and make the
#   new sheet the active tab, matching the final view state of the file

$wb = $excel.ActiveWorkbook

$hargeisa = $wb.Worksheets.Item("hargeisa")
$nigeria  = $wb.Worksheets.Item("nigeria")

# ---------------------------------------------------------------------
# New "sudan" sheet, inserted after "nigeria" (the last existing sheet)
# ---------------------------------------------------------------------
$sudan = $wb.Worksheets.Add($null, $nigeria)
$sudan.Name = "sudan"

$headers = @("indicator", "variable", "label")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $sudan.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$rows = @(
    @(1.1, 'I1_sec_day', 'Feeling safe at day'),
    @(1.1, 'I1_sec_inc', 'Experience security incident'),
    @(1.1, 'I1_sec_rep', 'Report security incident'),
    @(2.1, 'I3_no_borrow', 'Borrowing for food'),
    @(2.2, 'I4_hous_ownership', 'Own house'),
    @(2.2, 'I4_hous_water', 'Improved water '),
    @(2.2, 'I4_hous_toilet', 'Improved sanitation'),
    @(2.3, 'I5_med_satis', 'Satisfied with health facilities'),
    @(2.4, 'I6_ever_school', 'Ever in school'),
    @(2.4, 'I6_educ_child', 'Child in school'),
    @(3.1, 'I7_job_unemploy', 'Unemployment'),
    @(3.2, 'I8_econ_account', 'Bank account'),
    @(3.2, 'I8_poor190', 'Below 1.9 USD Poverty Line'),
    @(3.2, 'I8_poor32', 'Below 3.2 USD Poverty Line'),
    @(4.1, 'I9_hlp_access', 'Access to compensation'),
    @(4.1, 'I9_hlp_doc', 'Documentation'),
    @(5.1, 'I10_doc_birth', 'Birth certificate')
)

$r = 2
foreach ($row in $rows) {
    $sudan.Cells.Item($r, 1).Value = $row[0]
    $sudan.Cells.Item($r, 2).Value = $row[1]
    $sudan.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# ---------------------------------------------------------------------
# View / selection clean-up
# ---------------------------------------------------------------------
$hargeisa.Activate()
$hargeisa.Range("D11").Select()

$nigeria.Activate()
$nigeria.Range("C22").Select()

$sudan.Activate()
$sudan.Range("B22").Select()
